$wb = $excel.ActiveWorkbook

# --- "About" sheet: add India:US cost-adjustment note + factor -------------
$about = $wb.Worksheets.Item("About")

$about.Range("A29").Value = "India:US cost adjustment"
$about.Range("A30").Value = "see ""scaling-factors.xlsx in the InputData folder for source information."
$about.Range("A31").Value = 0.50596615326007366

# --- "DACD-capex" sheet: scale the US capex figure by the India:US factor --
$capex = $wb.Worksheets.Item("DACD-capex")
$capex.Range("B2").Formula = "=Data!B10*About!A31"

# --- Restore view/selection state seen in the target workbook --------------
$capex.Activate()
[void]$capex.Range("B3").Select()

$about.Activate()
[void]$about.Range("A29:A31").Select()
